# The Title, Author and Abstract paragraphs each had their text split
# word-by-word across many runs (one <w:r> per word / space). The edit
# merges each of those paragraphs back down to a single run carrying the
# full text, with no change in the visible wording.
#
# Word's `Range.Text = "..."` setter normally collapses a multi-run range
# into one run, but this host treats a same-text write as a no-op, so a
# direct reassignment leaves the per-word run split untouched. And this
# host's `Range.InsertXML` inserts its payload at the end of the
# paragraph rather than truly replacing the addressed range, so calling
# it directly on the populated range just appends a duplicate run.
#
# Deleting the run content first (collapsing the paragraph to empty, but
# keeping its paragraph mark/properties) and THEN inserting the desired
# single-run XML sidesteps both quirks and reproduces the exact
# xml:space="preserve" run markup used elsewhere in this document.

$d = $word.ActiveDocument

function Merge-ParagraphRuns($paragraph, [string]$styleId, [string]$text) {
    $escapedText = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

    $r = $paragraph.Range
    $r.MoveEnd(1, -1)
    $r.Delete()

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:pPr><w:pStyle w:val="' + $styleId + '"/></w:pPr>' +
           '<w:r><w:t xml:space="preserve">' + $escapedText + '</w:t></w:r></w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $r2 = $paragraph.Range
    $r2.MoveEnd(1, -1)
    $r2.InsertXML($xml)
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $styleName = $p.Style.NameLocal

    if ($styleName -eq "Title") {
        Merge-ParagraphRuns $p "Title" "Answers: Rationalizing the denominator"
    } elseif ($styleName -eq "Author") {
        Merge-ParagraphRuns $p "Author" "Maximilian Volmar"
    } elseif ($styleName -eq "Abstract") {
        Merge-ParagraphRuns $p "Abstract" "Answers to questions relating to the guide on rationalizing the denominator."
    }
}
